# Add November 2021 data (through 11-01), rolling the prior "Total" row
# into the new November row and appending a fresh Total row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet tab to reflect the new "through" date.
$ws.Name = "Through 2021-11-01"

# 2. The October row label loses its "(through 10-31)" suffix now that
#    the month is complete.
$ws.Range("A11").Value = "October"

# 3. Column A needs to be a bit wider for "November (through 11-01)".
$ws.Columns.Item(1).ColumnWidth = 23.8

# 4. Row 12 (previously the running "Total" row) becomes the November row.
$ws.Range("A12").Value = "November (through 11-01)"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 9
$ws.Range("H12").Value = 5

# 5. Append a new Total row at 13, carrying the same header style as the
#    other label cells in column A (copy A12's format before it is
#    overwritten below, so A13 keeps the bold/border/alignment formatting).
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A13").Value = "Total"
$ws.Range("B13").Value = 259
$ws.Range("C13").Value = 489
$ws.Range("D13").Value = 715
$ws.Range("E13").Value = 619
$ws.Range("F13").Value = 483
$ws.Range("G13").Value = 1066
$ws.Range("H13").Value = 1449
